# The document has several "<id>...</id>" markers that are each split across
# three runs: a Courier-New-styled "<id>" run, a plain black "p084r_aN" run,
# and a Courier-New-styled "</id>" run. The edit collapses each triplet into
# a single run (keeping the Courier New formatting of the surrounding runs)
# whose text is "<id>p084r_N</id>" (dropping the "a" from the id and merging
# the three runs/texts into one).

$d = $word.ActiveDocument

$ids = @("p084r_a1", "p084r_a3", "p084r_a4", "p084r_a5")

foreach ($oldId in $ids) {
    $newId = $oldId -replace "_a", "_"
    $oldTag = "<id>" + $oldId + "</id>"
    $newTag = "<id>" + $newId + "</id>"

    $r = $d.Content
    $found = $r.Find.Execute($oldTag, $false, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if ($found) {
        # Re-assigning the whole matched range's Text merges the three
        # underlying runs into a single run, inheriting the formatting
        # (Courier New, color 7f6000, sz/szCs 18) of the first run in the
        # match (the "<id>" run).
        $r.Text = $newTag
    }
}
